$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Update existing sheet "Transmittals_New": change the Action-Level2
#    value for the Single User (row 2) and Multi User (row 3) scenarios from
#    "Comments for Request for Information" to "Submission", fill in the
#    document-related columns (G:K) for row 3 (Multi User) to match row 2,
#    and append two new rows (4 and 5) describing the
#    Submission -> Forward -> Submission flow for Single User / Multi User.
# ---------------------------------------------------------------------------

$ws1.Range("M2").Value = "Submission"

$ws1.Range("G3").Value = "Document Register"
$ws1.Range("H3").Value = "Test 1 ta.docx"
$ws1.Range("I3").Value = "Document Register"
$ws1.Range("J3").Value = "Test 1 ta.docx"
$ws1.Range("K3").Value = "BrowseDocument.docx"
$ws1.Range("M3").Value = "Submission"

# Row 4 - Single User : Forward then Submission
$ws1.Range("A4").Value = "AutoTestAdmin"
$ws1.Range("C4").Value = "New Transmittal from Automation"
$ws1.Range("D4").Value = "UnTick"
$ws1.Range("E4").Value = "Correspondence"
$ws1.Range("F4").Value = "Request for Information"
$ws1.Range("L4").Value = "Message for New transmittal"
$ws1.Range("M4").Value = "Forward"
$ws1.Range("N4").Value = "AutoTestUser"
$ws1.Range("O4").Value = "Submission"

# Row 5 - Multi User : Forward then Submission
$ws1.Range("A5").Value = "AutoTestAdmin@@AutoTestUser"
$ws1.Range("C5").Value = "New Transmittal from Automation"
$ws1.Range("D5").Value = "UnTick"
$ws1.Range("E5").Value = "Correspondence"
$ws1.Range("F5").Value = "Request for Information"
$ws1.Range("L5").Value = "Message for New transmittal"
$ws1.Range("M5").Value = "Forward"
$ws1.Range("N5").Value = "SPInstall"
$ws1.Range("O5").Value = "Submission"

# ---------------------------------------------------------------------------
# 2) Add a new worksheet "Transmittals_New_ActionRequired" right after
#    "Transmittals_New". Build it off a copy of the first sheet so it keeps
#    the same header row, column widths and styles, then trim it down to
#    just the header row plus a single data row (the Single User scenario).
# ---------------------------------------------------------------------------

$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Transmittals_New_ActionRequired"

# Keep only header (row 1) + row 2; drop the old rows 3-5 that came from the copy
$ws2.Rows("3:5").Delete()

# Row 2 on the new sheet represents Submission for the Single User scenario
$ws2.Range("M2").Value = "Submission"

$ws1.Activate()
